$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename item in row 67: "Create table from CSV file" -> "Import table from CSV file"
$ws.Range("B67").Value = "Import table from CSV file"

# --- Add four new roadmap rows (72-75) ---
$ws.Rows.Item(72).RowHeight = 30

$ws.Cells.Item(72, 1).Value = 72
$ws.Cells.Item(72, 2).Value = "Create complementary block to link a post containing details from a row in a table back to the post with the table"
$ws.Cells.Item(72, 3).Value = "Feature"
$ws.Cells.Item(72, 4).Value = "All"
$ws.Cells.Item(72, 5).Value = "Summary"
$ws.Cells.Item(72, 6).Value = "Low"
$ws.Cells.Item(72, 7).Value = "Roadmap"
$ws.Cells.Item(72, 10).Value = "Roadmap"

$ws.Cells.Item(73, 1).Value = 73
$ws.Cells.Item(73, 2).Value = "Export table (data only)"
$ws.Cells.Item(73, 3).Value = "Feature"
$ws.Cells.Item(73, 4).Value = "All"
$ws.Cells.Item(73, 5).Value = "Summary"
$ws.Cells.Item(73, 6).Value = "Low"
$ws.Cells.Item(73, 7).Value = "Roadmap"
$ws.Cells.Item(73, 10).Value = "Roadmap"

$ws.Cells.Item(74, 1).Value = 74
$ws.Cells.Item(74, 2).Value = "Export table (full object)"
$ws.Cells.Item(74, 3).Value = "Feature"
$ws.Cells.Item(74, 4).Value = "All"
$ws.Cells.Item(74, 5).Value = "Summary"
$ws.Cells.Item(74, 6).Value = "Low"
$ws.Cells.Item(74, 7).Value = "Roadmap"
$ws.Cells.Item(74, 10).Value = "Roadmap"

$ws.Cells.Item(75, 1).Value = 75
$ws.Cells.Item(75, 2).Value = "Export all tables"
$ws.Cells.Item(75, 3).Value = "Feature"
$ws.Cells.Item(75, 4).Value = "All"
$ws.Cells.Item(75, 5).Value = "Summary"
$ws.Cells.Item(75, 6).Value = "Low"
$ws.Cells.Item(75, 7).Value = "Roadmap"
$ws.Cells.Item(75, 10).Value = "Roadmap"

# --- Extend AutoFilter range to include new row 72, add "Roadmap" to the Status filter list ---
$ws.AutoFilterMode = $false
$rng = $ws.Range("B1:K72")
$rng.AutoFilter(9, @("In Process", "Open", "Roadmap", "Testing"), 7)

# --- Update the _FilterDatabase defined name to match the new autofilter range ---
$fdName = $wb.Names.Item("Sheet1!_FilterDatabase")
$fdName.RefersTo = "=Sheet1!`$B`$1:`$K`$72"

# --- Update selection to reflect where the user last clicked ---
$ws.Activate()
$ws.Range("K70").Select()
